# Handback status report: refresh timestamps / priority as part of
# "Generate Report for Handback".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 3bd6b494... row (row 2) and the c3199b73... row (row 4) both shared the
# same timestamp string; update both occurrences.
$wsOverview.Range("G2").Value = "2016-08-29 20:26:18"
$wsOverview.Range("G4").Value = "2016-08-29 20:26:18"

# zh-cn sheet: Priority column (E) changes from "ht" to "mt" for both rows
# that used it.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback
# DateTime (K) timestamps.
$wsZhCn.Range("H2").Value = "2016-08-29 20:26:13"
$wsZhCn.Range("H4").Value = "2016-08-29 20:26:13"
$wsZhCn.Range("K2").Value = "2016-08-29 20:26:37"
$wsZhCn.Range("K4").Value = "2016-08-29 20:26:37"

# de-de sheet: Correspond Handoff Datetime (H) mirrors the Overview sheet's
# "Latest HO Xliff Generate Date" string.
$wsDeDe.Range("H2").Value = "2016-08-29 20:26:18"
$wsDeDe.Range("H4").Value = "2016-08-29 20:26:18"

# de-de sheet: Priority column (E) changes from "ht" to "mt" for both rows.
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# de-de sheet: Correspond Handback DateTime (K).
$wsDeDe.Range("K2").Value = "2016-08-29 20:26:44"
$wsDeDe.Range("K4").Value = "2016-08-29 20:26:44"
